## feat: add 2022-Q3 data
##
## - Inserts a new "2022-Q3" worksheet (between "总计" and "2022-Q2"),
##   duplicated from "2022-Q2" so it keeps identical layout/styling,
##   then overwritten with the Q3 fund figures.
## - Updates the "总计" summary sheet: the newest quarter (Q3) becomes
##   row 2, with the former row 2 (Q2) and row 3 (Q1) figures shifting
##   down to rows 3 and 4 respectively.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q3" worksheet right before "2022-Q2" by
# duplicating "2022-Q2" (this carries over sheetPr/pageMargins/styles
# exactly, instead of starting from a blank default sheet).
# ---------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($wsQ2)
$wsQ3 = $wb.Worksheets.Item("2022-Q2 (2)")
$wsQ3.Name = "2022-Q3"

# Replace the single fund row with the new quarter's figures.
$wsQ3.Range("C2").Value = "上投摩根亚太优势混合（QDII）"

# D2:G2 hold numbers-as-text in this workbook; force text storage so
# "22.76" etc. are written verbatim (no trailing float noise), then
# drop the temporary number format so the cells stay unstyled like
# their neighbours.
$wsQ3.Range("D2:G2").NumberFormat = "@"
$wsQ3.Range("D2").Value = "22.76"
$wsQ3.Range("E2").Value = "90.69"
$wsQ3.Range("F2").Value = "2.74"
$wsQ3.Range("G2").Value = "0.6236"
$wsQ3.Range("H2").Value = 6
$wsQ3.Range("D2:G2").ClearFormats()

# ---------------------------------------------------------------------
# Step 2: update the "总计" (totals) sheet so the quarters stay in
# descending order: Q3 (new), Q2, Q1.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Append a new row 4 carrying what used to be the 2022-Q1 entry.
# Copy A3's formatting (bold/border style) down to A4 first so the new
# row matches the look of the existing index column.
$wsTotal.Range("A3").Copy($wsTotal.Range("A4"))
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2022-Q1"
$wsTotal.Range("C4").Value = 1
$wsTotal.Range("D4").Value = 0.49

# Former row 2 (2022-Q2) now lives in row 3.
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("D3").Value = 0.55

# Row 2 becomes the brand-new 2022-Q3 entry.
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("D2").Value = 0.62

Write-Host "2022-Q3 data added"
